# Delete the entire first data row (spreadsheet row 2), shifting all
# subsequent rows up by one. This corresponds to removing the "530"
# record (id ad82673f-3f9f-4e02-a65d-a97c60514f2e) from the Notion export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
